$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.227.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "'2.520.68"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.05%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'541.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").Value = "'144.76"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").Value = "'2.543.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.36%  "

$ws.Range("D10").Value = "'0.100"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.35%  "

$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("E12").Value = "  +4.16%  "

$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").Value = "'2.966.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.52%  "

$ws.Range("D15").Value = "'23.73"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.12%  "

$ws.Range("D16").Value = "'59.158.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.40%  "

$ws.Range("E17").Value = "  +2.03%  "

$ws.Range("D18").Value = "'2.532.97"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.12%  "

$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").Value = "'4.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.56%  "

$ws.Range("D21").Value = "'324.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.47%  "

$ws.Range("E22").Value = "  +3.41%  "

$ws.Range("D23").Value = "'5.84"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.25%  "

$ws.Range("D24").Value = "'62.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.65%  "

$ws.Range("D25").Value = "'0.432"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.95%  "

$ws.Range("E26").Value = "  +2.88%  "

$ws.Range("E27").Value = "  +1.85%  "

$ws.Range("D28").Value = "'8.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.70%  "

$ws.Range("D29").Value = "'0.0₃0781"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").Value = "'1.82"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.22%  "

$ws.Range("D31").Value = "'6.70"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.30%  "

$ws.Range("E32").Value = "  -3.84%  "

$ws.Range("E33").Value = "  +7.13%  "

$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").Value = "'157.72"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.60%  "

$ws.Range("D36").Value = "'18.71"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.78%  "

$ws.Range("D37").Value = "'4.38"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("E38").Value = "  -5.40%  "

$ws.Range("E39").Value = "  -3.24%  "

$ws.Range("D40").Value = "'36.94"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'3.71"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'295.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.72%  "

$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("E44").Value = "  -0.20%  "

$ws.Range("D45").Value = "'0.601"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("D46").Value = "'10.81"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.80%  "

$ws.Range("E47").Value = "  -0.21%  "

$ws.Range("D48").Value = "'18.79"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.23%  "

$ws.Range("D49").Value = "'122.19"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.36%  "

$ws.Range("D50").Value = "'0.0515"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.25%  "

$ws.Range("E51").Value = "  -0.68%  "
